$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.319.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +9.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.760.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9959"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9930"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3759"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.89"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3425"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.200"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07588"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9949"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.434"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.104"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.756.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001099"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.75%  "
$ws.Range("E18").Value = "  +2.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9935"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.59%  "
$ws.Range("E22").Value = "  +6.00%  "
$ws.Range("B23").Value = "WrappedBTC"
$ws.Range("C23").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.294.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.40%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.447"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.504"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +26.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.442"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.956.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "133.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.120"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.107"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08646"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("E35").Value = "  +1.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "13.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.496"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02370"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.15%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06341"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.58%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2198"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.579"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6487"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.229"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9936"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("E46").Value = "  +9.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.937"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.101"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "130.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07271"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.66%  "
